$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "306.60"
Set-TextValue $ws.Range("E2") "-4.71%"
Set-TextValue $ws.Range("G2") "14"
Set-TextValue $ws.Range("D3") "38.97"
Set-TextValue $ws.Range("E3") "-8.72%"
Set-TextValue $ws.Range("G3") "14"
Set-TextValue $ws.Range("D4") "5.096"
Set-TextValue $ws.Range("E4") "-2.33%"
Set-TextValue $ws.Range("G4") "14"
Set-TextValue $ws.Range("D5") "0.07670"
Set-TextValue $ws.Range("E5") "-6.25%"
Set-TextValue $ws.Range("G5") "14"
Set-TextValue $ws.Range("D6") "4.251"
Set-TextValue $ws.Range("E6") "-1.49%"
Set-TextValue $ws.Range("G6") "14"
Set-TextValue $ws.Range("D7") "1.612"
Set-TextValue $ws.Range("E7") "-9.62%"
Set-TextValue $ws.Range("G7") "14"
Set-TextValue $ws.Range("D8") "0.9144"
Set-TextValue $ws.Range("E8") "-3.60%"
Set-TextValue $ws.Range("G8") "14"
Set-TextValue $ws.Range("D9") "0.1014"
Set-TextValue $ws.Range("E9") "-9.31%"
Set-TextValue $ws.Range("G9") "14"
Set-TextValue $ws.Range("D10") "0.1742"
Set-TextValue $ws.Range("E10") "-7.91%"
Set-TextValue $ws.Range("G10") "14"
Set-TextValue $ws.Range("D11") "0.09005"
Set-TextValue $ws.Range("E11") "-4.68%"
Set-TextValue $ws.Range("G11") "14"
Set-TextValue $ws.Range("D12") "0.04441"
Set-TextValue $ws.Range("E12") "-3.58%"
Set-TextValue $ws.Range("G12") "14"
Set-TextValue $ws.Range("D13") "0.1055"
Set-TextValue $ws.Range("E13") "-0.29%"
Set-TextValue $ws.Range("G13") "14"
Set-TextValue $ws.Range("D14") "0.001265"
Set-TextValue $ws.Range("E14") "-2.49%"
Set-TextValue $ws.Range("G14") "14"
Set-TextValue $ws.Range("D15") "0.005813"
Set-TextValue $ws.Range("E15") "0.24%"
Set-TextValue $ws.Range("G15") "14"
Set-TextValue $ws.Range("E16") "2,416.83%"
Set-TextValue $ws.Range("G16") "14"
Set-TextValue $ws.Range("D17") "3.360"
Set-TextValue $ws.Range("E17") "-0.02%"
Set-TextValue $ws.Range("G17") "14"
Set-TextValue $ws.Range("G18") "14"
Set-TextValue $ws.Range("D19") "0.3307"
Set-TextValue $ws.Range("E19") "-1.83%"
Set-TextValue $ws.Range("G19") "14"
Set-TextValue $ws.Range("D20") "7.031"
Set-TextValue $ws.Range("E20") "-5.71%"
Set-TextValue $ws.Range("G20") "14"
Set-TextValue $ws.Range("D21") "0.1348"
Set-TextValue $ws.Range("E21") "-2.86%"
Set-TextValue $ws.Range("G21") "14"
Set-TextValue $ws.Range("D22") "0.2817"
Set-TextValue $ws.Range("E22") "10.55%"
Set-TextValue $ws.Range("G22") "14"
Set-TextValue $ws.Range("D23") "0.04133"
Set-TextValue $ws.Range("E23") "-0.24%"
Set-TextValue $ws.Range("G23") "14"
Set-TextValue $ws.Range("D24") "0.001206"
Set-TextValue $ws.Range("E24") "-3.37%"
Set-TextValue $ws.Range("G24") "14"
Set-TextValue $ws.Range("D25") "0.004108"
Set-TextValue $ws.Range("E25") "-4.16%"
Set-TextValue $ws.Range("G25") "14"
Set-TextValue $ws.Range("D26") "0.0001302"
Set-TextValue $ws.Range("E26") "6.80%"
Set-TextValue $ws.Range("G26") "14"
Set-TextValue $ws.Range("G27") "14"
Set-TextValue $ws.Range("G28") "14"
Set-TextValue $ws.Range("G29") "14"
Set-TextValue $ws.Range("G30") "14"
Set-TextValue $ws.Range("G31") "14"
Set-TextValue $ws.Range("G32") "14"
Set-TextValue $ws.Range("G33") "14"
Set-TextValue $ws.Range("G34") "14"
Set-TextValue $ws.Range("G35") "14"
Set-TextValue $ws.Range("G36") "14"
Set-TextValue $ws.Range("G37") "14"
Set-TextValue $ws.Range("D38") "0.02435"
Set-TextValue $ws.Range("E38") "-9.12%"
Set-TextValue $ws.Range("G38") "14"
Set-TextValue $ws.Range("D39") "0.05177"
Set-TextValue $ws.Range("E39") "-7.42%"
Set-TextValue $ws.Range("G39") "14"
Set-TextValue $ws.Range("D40") "0.007914"
Set-TextValue $ws.Range("E40") "-2.68%"
Set-TextValue $ws.Range("G40") "14"
Set-TextValue $ws.Range("E41") "-6.07%"
Set-TextValue $ws.Range("G41") "14"
Set-TextValue $ws.Range("D42") "0.007134"
Set-TextValue $ws.Range("E42") "9.10%"
Set-TextValue $ws.Range("G42") "14"
Set-TextValue $ws.Range("D43") "0.001952"
Set-TextValue $ws.Range("E43") "-4.23%"
Set-TextValue $ws.Range("G43") "14"
Set-TextValue $ws.Range("D44") "0.008400"
Set-TextValue $ws.Range("E44") "9.85%"
Set-TextValue $ws.Range("G44") "14"
Set-TextValue $ws.Range("D45") "0.3329"
Set-TextValue $ws.Range("E45") "3.71%"
Set-TextValue $ws.Range("G45") "14"
Set-TextValue $ws.Range("D46") "0.00006419"
Set-TextValue $ws.Range("E46") "-4.99%"
Set-TextValue $ws.Range("G46") "14"
Set-TextValue $ws.Range("E47") "0.21%"
Set-TextValue $ws.Range("G47") "14"
Set-TextValue $ws.Range("D48") "0.003004"
Set-TextValue $ws.Range("E48") "-26.66%"
Set-TextValue $ws.Range("G48") "14"
Set-TextValue $ws.Range("D49") "0.004307"
Set-TextValue $ws.Range("E49") "39.90%"
Set-TextValue $ws.Range("G49") "14"
Set-TextValue $ws.Range("D50") "0.00002103"
Set-TextValue $ws.Range("E50") "0.21%"
Set-TextValue $ws.Range("G50") "14"
Set-TextValue $ws.Range("D51") "0.0002003"
Set-TextValue $ws.Range("E51") "0.21%"
Set-TextValue $ws.Range("G51") "14"
